$d = $word.ActiveDocument

# --- Paragraph 2: date line -> replace text (collapses the 3-run/proofErr mess into one run) ---
$d.Content.Find.Execute("2 October, 2019", $false, $false, $false, $false, $false, $true, 1, $false, "23 December 2019", 2) | Out-Null

# --- Paragraph 3: "HON 441 & 442" -> "Honor's Program" ---
$d.Content.Find.Execute("HON 441 & 442", $false, $false, $false, $false, $false, $true, 1, $false, "Honor's Program", 2) | Out-Null

# --- Apply Times New Roman to paragraphs 1-3 (non-empty, so direct formatting works) ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Name = "Times New Roman"
$p1b = $d.Paragraphs.Item(1)
$p1b.Range.Font.NameBi = "Times New Roman"

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Name = "Times New Roman"
$p2b = $d.Paragraphs.Item(2)
$p2b.Range.Font.NameBi = "Times New Roman"

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Font.Name = "Times New Roman"
$p3b = $d.Paragraphs.Item(3)
$p3b.Range.Font.NameBi = "Times New Roman"

# --- Paragraph 4: bookmark paragraph -> add bold title-start run before the bookmark,
#     drop right alignment, and make the whole paragraph (incl. mark) bold Times New Roman ---
$p4 = $d.Paragraphs.Item(4)
$insertPoint = $d.Range($p4.Range.Start, $p4.Range.Start)
$insertPoint.InsertBefore("Industry 4.0 and the ")

$p4b = $d.Paragraphs.Item(4)
$p4b.Format.Alignment = 0
$p4b.Range.Font.Name = "Times New Roman"
$p4c = $d.Paragraphs.Item(4)
$p4c.Range.Font.NameBi = "Times New Roman"
$p4d = $d.Paragraphs.Item(4)
$p4d.Range.Bold = 1
$p4e = $d.Paragraphs.Item(4)
$p4e.Range.BoldBi = 1

# --- Two trailing empty paragraphs (jc=right, Times New Roman), replacing the single
#     trailing empty paragraph that existed before. Use a temporary placeholder character
#     so paragraph-mark-only formatting (no visible run) can be applied, then strip it. ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertBefore("X")
$p5b = $d.Paragraphs.Item(5)
$p5b.Range.Font.Name = "Times New Roman"
$p5c = $d.Paragraphs.Item(5)
$p5c.Range.Font.NameBi = "Times New Roman"
$p5d = $d.Paragraphs.Item(5)
$p5d.Range.Collapse(0)
$p5d.Range.InsertParagraphAfter()

$xRange = $d.Range($p5.Range.Start, $p5.Range.Start + 1)
$xRange.Text = ""

$p6 = $d.Paragraphs.Item(6)
$yPoint = $d.Range($p6.Range.Start, $p6.Range.Start)
$yPoint.InsertBefore("Y")
$yRange = $d.Range($p6.Range.Start, $p6.Range.Start + 1)
$yRange.Text = ""
